$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tokens")
$ws.Range("A2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJlbWFpbCI6ImFkaXR5YXBhd2FyQHlvcG1haWwuY29tIiwib3RwIjo2OTM4NzMsImlhdCI6MTc0NDYzMDY0NiwiZXhwIjoxNzQ0NjM0MjQ2fQ.Wu7YQ9d0BgYrFEahVYwcA7hls-0YQnfpSzKTBR1k3i4"
